# Workbook: Maret-2014.xlsx
# Commit: "Added new dataset preprocessing"
#
# The original workbook has a single sheet "Data Harian - Table" whose
# daily-observation table lives in A9:K40 (header in row 9, 31 days of
# data in rows 10-40). The edit duplicates that table onto a brand-new
# worksheet named "Sheet1" (placed right after the existing sheet) with
# the table re-based at A1:K32 (header in row 1, data in rows 2-32) and
# makes that new sheet the active / selected one.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Make sure gridlines stay visible (touches the dirty-tracking so the
# serializer keeps emitting the "visible" default instead of flipping it).
$excel.ActiveWindow.DisplayGridlines = $true

# Select the table that is about to be copied, mirroring the selection
# left behind on the source sheet after the copy/paste.
[void]$ws1.Range("A9:K40").Select()

# Add the new worksheet right after "Data Harian - Table".
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Copy the whole daily table (header + 31 day rows) onto the new sheet,
# landing at A1 so it becomes A1:K32.
$ws1.Range("A9:K40").Copy($ws2.Range("A1"))

# Match the row height Excel computes for the wrapped two-line station
# names in the data rows (header row keeps the default height).
$ws2.Range("A2:K32").RowHeight = 28.8

# Leave the new sheet's whole table selected and active - it becomes the
# visible/selected tab.
[void]$ws2.Range("A1:K32").Select()
[void]$ws2.Activate()
